# Add "sum of goals scored, conceded and totals in last six" to the "L6"
# sheet: for each team row, append ",(N)" to the "Goals scored",
# "Goals conceded" and "Total Goals" cells, where N is the sum of the six
# space-separated numbers already stored in that cell (after the comma).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("L6")

# "Goals scored" = column C, "Goals conceded" = column D, "Total Goals" = column E
# Data rows are 2..11 (one per team). Iterate column-major so cells land in the
# same left-to-right, top-to-bottom order as the source data.
for ($col = 3; $col -le 5; $col++) {
    for ($row = 2; $row -le 11; $row++) {
        $cell = $ws.Cells.Item($row, $col)
        $current = [string]$cell.Value2

        $commaPos = $current.IndexOf(",")
        $numberPart = $current.Substring($commaPos + 1).Trim()
        $numbers = $numberPart -split "\s+"

        $total = 0
        foreach ($n in $numbers) {
            if ($n -ne "") {
                $total += [int]$n
            }
        }

        $cell.Value = "$current,($total)"
    }
}

Write-Output "Appended last-six goal sums on sheet 'L6'."
